$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save" (reuse G1's header style/format)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells H2:H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
